{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\n// 1) Title paragraph \"Group 10\" -> Heading 1, centered.\nconst titlePara = paras.items[0];\ntitlePara.style = \"Heading 1\";\ntitlePara.alignment = Word.Alignment.centered;\n\n// 2) Subtitle paragraph \"Homework 2\" -> Heading 1, centered.\nconst subtitlePara = paras.items[1];\nsubtitlePara.style = \"Heading 1\";\nsubtitlePara.alignment = Word.Alignment.centered;\n\nawait context.sync();\n\n// 3) Mention that Celery-Flower is now available, appended as a new\n//    sentence/run at the end of the \"installing and running\" paragraph.\nconst installPara = paras.items[4];\nconst newRange = installPara.insertText(\n  \"Celery-Flower is accessible on port 5555.\",\n  Word.InsertLocation.end\n);\n// Match the language formatting (\"en-US\") already used throughout the\n// document so the new run's rPr mirrors its neighbours.\nnewRange.languageId = \"en-US\";\nawait context.sync();\n\n// 4) Tighten up the wording of the celery-tests paragraph: drop the\n//    \"(during the actual test)\" aside and simplify \"most if not all\" to\n//    \"all\".\nconst oldSpan =\n  \"time-consuming (during the actual test) and finnicky, most if not all celery\";\nconst newSpan = \"time-consuming and finnicky, all celery\";\nconst results = body.search(oldSpan, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(newSpan, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Title paragraph \"Group 10\" -> Heading 1, centered.\n$titlePara = $d.Paragraphs(1)\n$titlePara.Range.Style = \"Heading 1\"\n$titlePara.Range.ParagraphFormat.Alignment = 1   # wdAlignParagraphCenter\n\n# 2) Subtitle paragraph \"Homework 2\" -> Heading 1, centered.\n$subtitlePara = $d.Paragraphs(2)\n$subtitlePara.Range.Style = \"Heading 1\"\n$subtitlePara.Range.ParagraphFormat.Alignment = 1   # wdAlignParagraphCenter\n\n# 3) Mention that Celery-Flower is now available, appended as a new\n#    sentence/run at the end of the \"installing and running\" paragraph\n#    (paragraph 5, right before its trailing paragraph mark).\n$installPara = $d.Paragraphs(5)\n$pRange = $installPara.Range\n$insertionPoint = $d.Range($pRange.Start, $pRange.End - 1)\n$insertionPoint.Collapse(0)   # wdCollapseEnd\n$insertionPoint.InsertAfter(\"Celery-Flower is accessible on port 5555.\")\n# Match the language formatting (\"en-US\") already used throughout the\n# document so the new run's rPr mirrors its neighbours.\n$insertionPoint.LanguageID = \"en-US\"\n\n# 4) Tighten up the wording of the celery-tests paragraph: drop the\n#    \"(during the actual test)\" aside and simplify \"most if not all\" to\n#    \"all\".\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n  \"time-consuming (during the actual test) and finnicky, most if not all celery\",\n  $true,\n  $false,\n  $false,\n  $false,\n  $false,\n  $true,\n  1,\n  $false,\n  \"time-consuming and finnicky, all celery\",\n  2\n)\n"}
